$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.365.31"
$ws.Range("E2").Value = "  +5.92%  "
$ws.Range("D3").Value = "3.554.13"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormatLocal = "@"
$ws.Range("D5").Value = "417.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormatLocal = "@"
$ws.Range("D6").Value = "130.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormatLocal = "@"
$ws.Range("D7").Value = "0.657"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.86%  "
$ws.Range("D8").Value = "3.547.65"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormatLocal = "@"
$ws.Range("D10").Value = "0.776"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.25%  "
$ws.Range("D11").NumberFormatLocal = "@"
$ws.Range("D11").Value = "0.179"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +27.38%  "
$ws.Range("D12").NumberFormatLocal = "@"
$ws.Range("D12").Value = "0.0000331"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +50.15%  "
$ws.Range("D13").NumberFormatLocal = "@"
$ws.Range("D13").Value = "43.12"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").NumberFormatLocal = "@"
$ws.Range("D14").Value = "10.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.22%  "
$ws.Range("D15").Value = "4.114.63"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormatLocal = "@"
$ws.Range("D17").Value = "20.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "3.591.48"
$ws.Range("E18").Value = "  +4.71%  "
$ws.Range("D19").NumberFormatLocal = "@"
$ws.Range("D19").Value = "1.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.93%  "
$ws.Range("D20").NumberFormatLocal = "@"
$ws.Range("D20").Value = "12.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").Value = "66.228.72"
$ws.Range("E21").Value = "  +5.84%  "
$ws.Range("D22").NumberFormatLocal = "@"
$ws.Range("D22").Value = "452.93"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").NumberFormatLocal = "@"
$ws.Range("D23").Value = "90.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormatLocal = "@"
$ws.Range("D24").Value = "3.20"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormatLocal = "@"
$ws.Range("D25").Value = "13.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").NumberFormatLocal = "@"
$ws.Range("D26").Value = "3.38"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormatLocal = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").NumberFormatLocal = "@"
$ws.Range("D28").Value = "34.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("D29").NumberFormatLocal = "@"
$ws.Range("D29").Value = "4.83"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormatLocal = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormatLocal = "@"
$ws.Range("D31").Value = "12.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.13%  "
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").NumberFormatLocal = "@"
$ws.Range("D33").Value = "7.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.32%  "
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("D35").NumberFormatLocal = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormatLocal = "@"
$ws.Range("D36").Value = "39.03"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0815"
$ws.Range("E37").Value = "  +45.65%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormatLocal = "@"
$ws.Range("D38").Value = "56.87"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormatLocal = "@"
$ws.Range("D39").Value = "0.0499"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormatLocal = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.56%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").NumberFormatLocal = "@"
$ws.Range("D43").Value = "3.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormatLocal = "@"
$ws.Range("D44").Value = "149.06"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("E47").Value = "  -5.13%  "
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").NumberFormatLocal = "@"
$ws.Range("D49").Value = "2.43"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").NumberFormatLocal = "@"
$ws.Range("D50").Value = "15.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.09%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormatLocal = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +10.34%  "
